$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.211271643638611
$ws.Range("B1").Value = 2.415652513504028
$ws.Range("C1").Value = 7.065585136413574
$ws.Range("D1").Value = 2.268170833587646
$ws.Range("E1").Value = 1.168520212173462
